# Rename the existing sheet and add a new "Ví dụ" (example) sheet that
# documents common mistakes, mirroring the student_info layout.

$wb = $excel.ActiveWorkbook

# 1) Rename Sheet1 -> student_info
$wsInfo = $wb.Worksheets.Item(1)
$wsInfo.Name = "student_info"

# 2) Move the (no longer active) selection on student_info to A22
$wsInfo.Range("A22").Select() | Out-Null

# 3) Add the new sheet right after student_info; it becomes the active tab
$wsExample = $wb.Worksheets.Add($null, $wsInfo)
$wsExample.Name = "Ví dụ"

# Column widths to match student_info-like layout
$wsExample.Columns("A:G").ColumnWidth = 15.44140625

# Header row
$headers = @("student_name", "username", "password", "course", "level", "lesson_num", "image_name")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsExample.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

# Row 2 sample data
$wsExample.Range("A2").Value = "Phạm Vy An"
$wsExample.Range("B2").Value = 84972235270
$wsExample.Range("C2").Value = "Teky@123"
$wsExample.Range("D2").Value = "Bé làm Game"
$wsExample.Range("E2").Value = 7
$wsExample.Range("F2").Value = 4
$wsExample.Range("G2").Value = "test.png"
$wsExample.Range("G2").HorizontalAlignment = -4152

# Row 3 sample data
$wsExample.Range("A3").Value = "Phan Quốc Hưng"
$wsExample.Range("B3").Value = 84937285555
$wsExample.Range("C3").Value = "Teky@123"
$wsExample.Range("D3").Value = "Bé làm Game"
$wsExample.Range("E3").Value = 7
$wsExample.Range("F3").Value = 5
$wsExample.Range("G3").Value = "test.png"
$wsExample.Range("G3").HorizontalAlignment = -4152

# Hyperlinks on the password column (mirrors student_info's C2/C3)
$wsExample.Hyperlinks.Add($wsExample.Range("C2"), "mailto:Teky@123") | Out-Null
$wsExample.Hyperlinks.Add($wsExample.Range("C3"), "mailto:Teky@123") | Out-Null
$wsExample.Range("C2").HorizontalAlignment = -4152
$wsExample.Range("C3").HorizontalAlignment = -4152

# Notes / common-mistakes section
$wsExample.Range("B12").Value = "Một số lỗi thường gặp"
$wsExample.Range("B12").Font.Bold = $true

# Note: the shared-string table records "image name" text before the
# "course name" text (matches original authoring order), even though the
# image-name note ends up one row below the course-name note.
$wsExample.Range("B14").NumberFormat = "@"
$wsExample.Range("B14").Value = "- Tên hình ảnh sản phẩm không giống với tên file hình trong thư mục, hoặc quên ghi đuôi .png/.jpg/.jpeg"

$wsExample.Range("B13").NumberFormat = "@"
$wsExample.Range("B13").Value = "- Tên khóa học ký tự viết hoa chưa giống"

$wsExample.Range("G6").Select() | Out-Null
